$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Anm"
$ws.Range("B12").Value = "Sakib"
